$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8169493941477981
$ws.Range("C2").Value = 0.2092783737809327
$ws.Range("E2").Value = 0.1175178547153592
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.00238735565926681
$ws.Range("M2").Value = 0.3392577148065001
$ws.Range("O2").Value = 1.400592669394598
$ws.Range("B3").Value = 0.7157122905342135
$ws.Range("C3").Value = 0.1866659086233255
$ws.Range("E3").Value = 0.1126404932284544
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.002390247812993586
$ws.Range("M3").Value = 0.3013087131993402
$ws.Range("O3").Value = 1.422772324373781
$ws.Range("B4").Value = 0.6533906370119666
$ws.Range("C4").Value = 0.1727166345719695
$ws.Range("E4").Value = 0.1097763865914061
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.002392116098657474
$ws.Range("M4").Value = 0.2780441076107536
$ws.Range("O4").Value = 1.438177764951178
$ws.Range("B5").Value = 0.6279546662132987
$ws.Range("C5").Value = 0.1670161627931463
$ws.Range("E5").Value = 0.1086417655814529
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.002392900769607682
$ws.Range("M5").Value = 0.2685727353118068
$ws.Range("O5").Value = 1.444902954433047
$ws.Range("B6").Value = 0.6237287010290515
$ws.Range("C6").Value = 0.1660686460082275
$ws.Range("E6").Value = 0.1084553178083709
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.002393032474630021
$ws.Range("M6").Value = 0.2670005745037329
$ws.Range("O6").Value = 1.446046626492759
$ws.Range("B7").Value = 0.6530477564214152
$ws.Range("C7").Value = 0.1726398204527015
$ws.Range("E7").Value = 0.1097609534202775
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.002392126586422775
$ws.Range("M7").Value = 0.2779163363016508
$ws.Range("O7").Value = 1.43826665460081
$ws.Range("B8").Value = 0.7820771138566442
$ws.Range("C8").Value = 0.2014953137271789
$ws.Range("E8").Value = 0.1158088385134022
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.002388333720724924
$ws.Range("M8").Value = 0.3261653841679362
$ws.Range("O8").Value = 1.407868100628207
$ws.Range("B9").Value = 1.033778125808453
$ws.Range("C9").Value = 0.2575516504853397
$ws.Range("E9").Value = 0.1287201521231793
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.00238162650200551
$ws.Range("M9").Value = 0.4210728691165002
$ws.Range("O9").Value = 1.362526962267779
$ws.Range("B10").Value = 1.217857016578876
$ws.Range("C10").Value = 0.298401128906761
$ws.Range("E10").Value = 0.1388689903986844
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.002377139450225376
$ws.Range("M10").Value = 0.4909915936955684
$ws.Range("O10").Value = 1.338040121227664
$ws.Range("B11").Value = 1.301409017770311
$ws.Range("C11").Value = 0.3169095383642286
$ws.Range("E11").Value = 0.1436346932775194
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.002375192887763537
$ws.Range("M11").Value = 0.5228439869388382
$ws.Range("O11").Value = 1.328843759633969
$ws.Range("B12").Value = 1.333020242912767
$ws.Range("C12").Value = 0.3239072431668717
$ws.Range("E12").Value = 0.1454611262697867
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.002374469306671813
$ws.Range("M12").Value = 0.5349123877010271
$ws.Range("O12").Value = 1.325642911820239
$ws.Range("B13").Value = 1.326213470400319
$ws.Range("C13").Value = 0.3224006583058099
$ws.Range("E13").Value = 0.1450667977200908
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.002374624541699649
$ws.Range("M13").Value = 0.5323129488946705
$ws.Range("O13").Value = 1.326319711143839
$ws.Range("B14").Value = 1.304010262910822
$ws.Range("C14").Value = 0.3174854667834097
$ws.Range("E14").Value = 0.1437845168311185
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.002375133087342158
$ws.Range("M14").Value = 0.5238367295732758
$ws.Range("O14").Value = 1.328574767162024
$ws.Range("B15").Value = 1.290406458988627
$ws.Range("C15").Value = 0.3144733234227601
$ws.Range("E15").Value = 0.1430019274825156
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.002375446347596347
$ws.Range("M15").Value = 0.5186456581261325
$ws.Range("O15").Value = 1.329992795051993
$ws.Range("B16").Value = 1.212392834273373
$ws.Range("C16").Value = 0.2971900367860769
$ws.Range("E16").Value = 0.1385605672657277
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.002377268561070587
$ws.Range("M16").Value = 0.4889108884531197
$ws.Range("O16").Value = 1.338680411120606
$ws.Range("B17").Value = 1.164485327832097
$ws.Range("C17").Value = 0.2865680355234588
$ws.Range("E17").Value = 0.1358743312947794
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.002378410617027018
$ws.Range("M17").Value = 0.4706813353829489
$ws.Range("O17").Value = 1.344509179930157
$ws.Range("B18").Value = 1.136912727861557
$ws.Range("C18").Value = 0.2804515759621609
$ws.Range("E18").Value = 0.1343432780673339
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.002379076407276104
$ws.Range("M18").Value = 0.4602005109381082
$ws.Range("O18").Value = 1.348044467757248
$ws.Range("B19").Value = 1.127574158799121
$ws.Range("C19").Value = 0.2783794624497773
$ws.Range("E19").Value = 0.1338272834950729
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.002379303364905898
$ws.Range("M19").Value = 0.4566526288373183
$ws.Range("O19").Value = 1.3492727720039
$ws.Range("B20").Value = 1.169586981167072
$ws.Range("C20").Value = 0.2876994891715867
$ws.Range("E20").Value = 0.1361588343864852
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.002378288121411316
$ws.Range("M20").Value = 0.4726214542596239
$ws.Range("O20").Value = 1.343869769007682
$ws.Range("B21").Value = 1.310532656160944
$ws.Range("C21").Value = 0.3189294797821844
$ws.Range("E21").Value = 0.1441605603535479
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.002374983348161656
$ws.Range("M21").Value = 0.5263262219632878
$ws.Range("O21").Value = 1.32790474176889
$ws.Range("B22").Value = 1.402484400018579
$ws.Range("C22").Value = 0.3392756426386825
$ws.Range("E22").Value = 0.1495171742291248
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.002372902378112298
$ws.Range("M22").Value = 0.5614638000946286
$ws.Range("O22").Value = 1.319113177534945
$ws.Range("B23").Value = 1.353423490627392
$ws.Range("C23").Value = 0.3284225248086159
$ws.Range("E23").Value = 0.1466465110719071
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.00237400583442919
$ws.Range("M23").Value = 0.5427067061326909
$ws.Range("O23").Value = 1.323654372325251
$ws.Range("B24").Value = 1.167280616850405
$ws.Range("C24").Value = 0.287187989259877
$ws.Range("E24").Value = 0.1360301691402341
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.002378343472990458
$ws.Range("M24").Value = 0.4717443278027531
$ws.Range("O24").Value = 1.344158272914882
$ws.Range("B25").Value = 0.9658322061433751
$ws.Range("C25").Value = 0.2424448945706388
$ws.Range("E25").Value = 0.1251124759512905
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.002383363249947042
$ws.Range("M25").Value = 0.395365415154231
$ws.Range("O25").Value = 1.373252201870187
